# ClueLayout_Highlighted.xlsx edit:
# "Changed walkway condition to make walkway W instead of A"
#
# Semantics: cells that currently hold "A" become "W" (the new walkway
# marker), cells that currently hold "W" become "O", and cells that hold
# "WD" become "OD" (same relabeling, applied to the door-adjacent variant).
# Every other cell value is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange

# Snapshot the (cell, newValue) pairs up front from the ORIGINAL values so
# that writing "W" into a former "A" cell can't be mistaken, later in the
# same pass, for an original "W" that should become "O".
$targets = @()
foreach ($cell in $used.Cells) {
    $v = $cell.Value2
    if ($v -eq "A") {
        $targets += ,@($cell, "W")
    } elseif ($v -eq "W") {
        $targets += ,@($cell, "O")
    } elseif ($v -eq "WD") {
        $targets += ,@($cell, "OD")
    }
}

foreach ($pair in $targets) {
    $pair[0].Value = $pair[1]
}

# Match the saved cursor/selection position recorded in the sheet.
$ws.Range("X26").Select() | Out-Null
